# Applies the "Se modifica arquitectura para presentar los reportes /
# Se comienza con el bundle de articulos" edit to Backlog.xlsx
# (Hoja1 = the backlog sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Add the new backlog row (row 101) -------------------------------
# New task string (becomes a new shared-string entry) with status
# "no comenzado" (column B), same pattern as all the other rows.
$ws.Range("A101").Value = "clase para simplificar el servidor de reportes"
$ws.Range("B101").Value = "no comenzado"

# --- 2. Hide the rows that are now filtered out --------------------------
# These tasks moved to "terminado"/out-of-view state and their rows become
# hidden (autoFilter hides them because they no longer match the
# "no comenzado" filter criteria).
$ws.Rows.Item(67).Hidden = $true
$ws.Rows.Item(87).Hidden = $true
$ws.Rows.Item(88).Hidden = $true
$ws.Rows.Item(90).Hidden = $true
$ws.Rows.Item(93).Hidden = $true

# --- 3. Re-apply the AutoFilter over the extended range -------------------
# Turn off the existing filter first, then reapply it across A1:C100 with
# the same "no comenzado" criteria on column B (field 2).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:C100").AutoFilter(2, "no comenzado", 7)

# --- 4. Extend the hidden _FilterDatabase defined name --------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$C`$100"
    }
}

# --- 5. Update the selection to the newly added cell -----------------------
$ws.Activate()
[void]$ws.Range("A101").Select()

Write-Host "edit applied"
